# Including NAMPT in POI
#
# The POI table (Sheet 1, A1:G10) is re-sorted by PPI_rank (col G) and a new
# row for NAMPT (P43490) is inserted, which shifts/re-derives GO_POI (D),
# POI_Go_found (E) and PPI_rank (G) for every data row. Rather than trying to
# reproduce the row-shuffling in place, every data cell A2:G11 is written
# directly with its final target value (the header row A1:G1 is unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - AADAT (Q8N5Z0)
$ws.Cells.Item(2,1).Value = 'Q8N5Z0'
$ws.Cells.Item(2,2).Value = 'aminoadipate aminotransferase [Source:HGNC Symbol;Acc:HGNC:17929]'
$ws.Cells.Item(2,3).Value = 'AADAT'
$ws.Cells.Item(2,4).Value = 0.0555555555555556
$ws.Cells.Item(2,5).Value = 'mitochondrion'
$ws.Cells.Item(2,6).Value = 'AADAT_HUMAN'
$ws.Cells.Item(2,7).Value = 0.0150975810802536

# Row 3 - VDAC1 (P21796)
$ws.Cells.Item(3,1).Value = 'P21796'
$ws.Cells.Item(3,2).Value = 'voltage dependent anion channel 1 [Source:HGNC Symbol;Acc:HGNC:12669]'
$ws.Cells.Item(3,3).Value = 'P21796'
$ws.Cells.Item(3,4).Value = 0.0526315789473684
$ws.Cells.Item(3,5).Value = 'membrane,mitochondrion'
$ws.Cells.Item(3,6).Value = 'VDAC1_HUMAN'
$ws.Cells.Item(3,7).Value = 6.09422906604821

# Row 4 - ISCU (Q9H1K1)
$ws.Cells.Item(4,1).Value = 'Q9H1K1'
$ws.Cells.Item(4,2).Value = 'iron-sulfur cluster assembly enzyme [Source:HGNC Symbol;Acc:HGNC:29882]'
$ws.Cells.Item(4,3).Value = 'ISCU_HUMAN'
$ws.Cells.Item(4,4).Value = 0.05
$ws.Cells.Item(4,5).Value = 'mitochondrion'
$ws.Cells.Item(4,6).Value = 'ISCU_HUMAN'
$ws.Cells.Item(4,7).Value = 2.90489511964124

# Row 5 - HXK2 (P52789)
$ws.Cells.Item(5,1).Value = 'P52789'
$ws.Cells.Item(5,2).Value = 'hexokinase 2 [Source:HGNC Symbol;Acc:HGNC:4923]'
$ws.Cells.Item(5,3).Value = 'P52789'
$ws.Cells.Item(5,4).Value = 0.0465116279069767
$ws.Cells.Item(5,5).Value = 'membrane,mitochondrion'
$ws.Cells.Item(5,6).Value = 'HXK2_HUMAN'
$ws.Cells.Item(5,7).Value = 3.30671353906589

# Row 6 - AAK1 (Q2M2I8)
$ws.Cells.Item(6,1).Value = 'Q2M2I8'
$ws.Cells.Item(6,2).Value = 'AP2 associated kinase 1 [Source:HGNC Symbol;Acc:HGNC:19679]'
$ws.Cells.Item(6,3).Value = 'AAK1'
$ws.Cells.Item(6,4).Value = 0.032258064516129
$ws.Cells.Item(6,5).Value = 'membrane'
$ws.Cells.Item(6,6).Value = 'AAK1_HUMAN'
$ws.Cells.Item(6,7).Value = 3.02365118426608

# Row 7 - NAMPT (P43490) - newly added
$ws.Cells.Item(7,1).Value = 'P43490'
$ws.Cells.Item(7,2).Value = 'nicotinamide phosphoribosyltransferase [Source:HGNC Symbol;Acc:HGNC:30092]'
$ws.Cells.Item(7,3).Value = 'NAMPT'
$ws.Cells.Item(7,4).Value = 0
$ws.Cells.Item(7,5).Value = ''
$ws.Cells.Item(7,6).Value = 'NAMPT_HUMAN'
$ws.Cells.Item(7,7).Value = 2.31334502757633

# Row 8 - PHF8 (Q9UPP1)
$ws.Cells.Item(8,1).Value = 'Q9UPP1'
$ws.Cells.Item(8,2).Value = 'PHD finger protein 8 [Source:HGNC Symbol;Acc:HGNC:20672]'
$ws.Cells.Item(8,3).Value = 'Q9UPP1'
$ws.Cells.Item(8,4).Value = 0
$ws.Cells.Item(8,5).Value = ''
$ws.Cells.Item(8,6).Value = 'PHF8_HUMAN'
$ws.Cells.Item(8,7).Value = 2.27543123935373

# Row 9 - SPCS (Q9HD40)
$ws.Cells.Item(9,1).Value = 'Q9HD40'
$ws.Cells.Item(9,2).Value = 'Sep (O-phosphoserine) tRNA:Sec (selenocysteine) tRNA synthase [Source:HGNC Symbol;Acc:HGNC:30605]'
$ws.Cells.Item(9,3).Value = 'SPCS_HUMAN'
$ws.Cells.Item(9,4).Value = 0
$ws.Cells.Item(9,5).Value = ''
$ws.Cells.Item(9,6).Value = 'SPCS_HUMAN'
$ws.Cells.Item(9,7).Value = 1.37959306612601

# Row 10 - AACS (Q86V21)
$ws.Cells.Item(10,1).Value = 'Q86V21'
$ws.Cells.Item(10,2).Value = 'acetoacetyl-CoA synthetase [Source:HGNC Symbol;Acc:HGNC:21298]'
$ws.Cells.Item(10,3).Value = 'AACS'
$ws.Cells.Item(10,4).Value = 0
$ws.Cells.Item(10,5).Value = ''
$ws.Cells.Item(10,6).Value = 'AACS_HUMAN'
$ws.Cells.Item(10,7).Value = 0.340125586489395

# Row 11 - SEPP1 (P49908)
$ws.Cells.Item(11,1).Value = 'P49908'
$ws.Cells.Item(11,2).Value = 'selenoprotein P [Source:HGNC Symbol;Acc:HGNC:10751]'
$ws.Cells.Item(11,3).Value = 'SEPP1_HUMAN'
$ws.Cells.Item(11,4).Value = 0
$ws.Cells.Item(11,5).Value = ''
$ws.Cells.Item(11,6).Value = 'SEPP1_HUMAN'
$ws.Cells.Item(11,7).Value = 0.308935676622387
